# Apply updated "dSF" (column F) values on the worksheet, per the data repull
# described in the commit message ("repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = -3
    4  = -1
    5  = 9
    6  = -4
    7  = -6
    9  = -5
    10 = 1
    11 = -1
    12 = -2
    13 = -1
    14 = -2
    16 = -4
    17 = -2
    18 = 2
    19 = 1
    21 = 3
    22 = -4
    24 = -4
    25 = -2
    26 = -2
    27 = -2
    28 = 2
    29 = 1
    30 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
